# "Dividend Testing - Error"
# Update a handful of tracked game-state values on Sheet1 (Balance.xlsx)
# and move the selection to reflect where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5  = Property Value  -> Player 3 (column C) goes from 0 to 10
$ws.Range("C5").Value = 10

# Row 10 = Net Worth -> Player 1 drops to 8, Player 2 jumps to 21,
#          Players 3-5 settle at 8 each
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = 21
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 8

# Row 11 = Turns Played -> every player now shows 8 turns played
$ws.Range("B11").Value = 8
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 8

# Row 14 = Debt Taken -> Player 3 takes on 2
$ws.Range("C14").Value = 2

# Row 15 = Debt value to be repaid -> Player 3 owes 2
$ws.Range("C15").Value = 2

# Reflect the cursor's final resting place on the sheet
$ws.Range("E11").Select() | Out-Null
